$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data row (row 2) values for the new Prepago test case ---
$ws.Range("A2").Value = "77513"
$ws.Range("B2").Value = "8757940"
$ws.Range("C2").Value = "14"
$ws.Range("J2").Value = "14/01/2022"

# --- Add new columns K and L with headers (row 1) ---
$ws.Range("K1").Value = "AccountingSourcePrepag"
$ws.Range("L1").Value = "AcountingNamePrepag"

# --- Add new values for the Prepago accounting source/name (row 2) ---
$ws.Range("K2").Value = """'PREPAG'"""
$ws.Range("L2").Value = """upper('Recaudo prepago')"""

# Ensure the new data cells (row 2) use the same "text" style as the rest of row 2
$ws.Range("K2").NumberFormat = "@"
$ws.Range("L2").NumberFormat = "@"

# --- Adjust column widths to better fit the new / resized columns ---
$ws.Columns.Item(3).ColumnWidth = 16.333333333333336
$ws.Columns.Item(10).ColumnWidth = 12.833333333333332
$ws.Columns.Item(11).ColumnWidth = 24.5
$ws.Columns.Item(12).ColumnWidth = 25.833333333333336

# --- Update the active selection / scroll position ---
$ws.Range("K2").Select()
